# Updated symbol list on Sun Jan 22 23:57:52 UTC 2023 with GitHub Actions
function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextCell $ws "D2" '302.68'
Set-TextCell $ws "E2" '0.71%'
# Row 3
Set-TextCell $ws "D3" '36.56'
Set-TextCell $ws "E3" '3.02%'
# Row 4
Set-TextCell $ws "D4" '4.986'
Set-TextCell $ws "E4" '-0.75%'
# Row 5
Set-TextCell $ws "D5" '0.07759'
# Row 6
Set-TextCell $ws "E6" '-3.45%'
# Row 7
Set-TextCell $ws "D7" '7.911'
Set-TextCell $ws "E7" '-1.40%'
# Row 8
Set-TextCell $ws "B8" 'MXToken'
Set-TextCell $ws "C8" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws "D8" '0.9218'
Set-TextCell $ws "E8" '-0.48%'
# Row 9
Set-TextCell $ws "B9" 'LiechtensteinCryptoassetsExchange'
Set-TextCell $ws "C9" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell $ws "D9" '0.09799'
Set-TextCell $ws "E9" '6.79%'
# Row 10
Set-TextCell $ws "B10" 'WazirX'
Set-TextCell $ws "C10" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell $ws "D10" '0.1868'
Set-TextCell $ws "E10" '2.40%'
# Row 11
Set-TextCell $ws "B11" 'MandalaExchangeToken'
Set-TextCell $ws "C11" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell $ws "D11" '0.08580'
Set-TextCell $ws "E11" '1.15%'
# Row 12
Set-TextCell $ws "B12" 'BitrueCoin'
Set-TextCell $ws "C12" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell $ws "D12" '0.03512'
Set-TextCell $ws "E12" '-1.41%'
# Row 13
Set-TextCell $ws "B13" 'BitMartToken'
Set-TextCell $ws "C13" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell $ws "D13" '0.09954'
Set-TextCell $ws "E13" '-0.08%'
# Row 14
Set-TextCell $ws "B14" 'BitForexToken'
Set-TextCell $ws "C14" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell $ws "D14" '0.001467'
Set-TextCell $ws "E14" '-0.80%'
# Row 15
Set-TextCell $ws "B15" 'TigerCash'
Set-TextCell $ws "C15" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws "D15" '0.005644'
Set-TextCell $ws "E15" '-1.99%'
# Row 16
Set-TextCell $ws "B16" 'LEO'
Set-TextCell $ws "C16" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws "D16" '3.464'
Set-TextCell $ws "E16" '-0.31%'
# Row 17
Set-TextCell $ws "B17" 'GateToken'
Set-TextCell $ws "C17" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell $ws "D17" '4.043'
Set-TextCell $ws "E17" '0.93%'
# Row 18
Set-TextCell $ws "D18" '2.311'
Set-TextCell $ws "E18" '5.83%'
# Row 19
Set-TextCell $ws "E19" '-1.57%'
# Row 20
Set-TextCell $ws "D20" '0.1343'
Set-TextCell $ws "E20" '1.41%'
# Row 21
Set-TextCell $ws "D21" '4.754'
Set-TextCell $ws "E21" '3.53%'
# Row 22
Set-TextCell $ws "D22" '0.2197'
Set-TextCell $ws "E22" '-2.09%'
# Row 23
Set-TextCell $ws "D23" '0.04576'
Set-TextCell $ws "E23" '-1.95%'
# Row 24
Set-TextCell $ws "D24" '0.005079'
Set-TextCell $ws "E24" '13.51%'
# Row 25
Set-TextCell $ws "D25" '0.001229'
Set-TextCell $ws "E25" '-0.74%'
# Row 26
Set-TextCell $ws "E26" '6.66%'
# Row 39
Set-TextCell $ws "D39" '0.01764'
Set-TextCell $ws "E39" '2.33%'
# Row 40
Set-TextCell $ws "D40" '0.04668'
Set-TextCell $ws "E40" '-0.12%'
# Row 41
Set-TextCell $ws "D41" '0.007455'
Set-TextCell $ws "E41" '-6.00%'
# Row 42
Set-TextCell $ws "E42" '-0.62%'
# Row 43
Set-TextCell $ws "D43" '0.007702'
Set-TextCell $ws "E43" '0.18%'
# Row 44
Set-TextCell $ws "E44" '-0.10%'
# Row 45
Set-TextCell $ws "E45" '16.42%'
# Row 46
Set-TextCell $ws "D46" '0.00006175'
Set-TextCell $ws "E46" '-0.94%'
# Row 47
Set-TextCell $ws "D47" '0.00000000749'
Set-TextCell $ws "E47" '-0.81%'
# Row 48
Set-TextCell $ws "D48" '0.0005793'
Set-TextCell $ws "E48" '-0.14%'
# Row 49
Set-TextCell $ws "D49" '38.16'
Set-TextCell $ws "E49" '1,040.09%'
# Row 50
Set-TextCell $ws "D50" '0.001998'
Set-TextCell $ws "E50" '-26.11%'
# Row 51
Set-TextCell $ws "D51" '0.00002098'
Set-TextCell $ws "E51" '-0.81%'
